$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AA2").Value = 0.00499322211208687
$ws.Range("AB2").Value = 0.001947497312168793
$ws.Range("AD2").Value = 0.01211462069384904
$ws.Range("AE2").Value = 0.004234256159415955
$ws.Range("AF2").Value = 0.02519078860322073
$ws.Range("AH2").Value = 0.0004139410213846818
$ws.Range("AI2").Value = 0.02926816432373484
$ws.Range("D2").Value = 0.2022278673415956
$ws.Range("E2").Value = 0.007487481361975501
$ws.Range("F2").Value = 0.3883437733773492
$ws.Range("H2").Value = 0.0299525876193413
$ws.Range("J2").Value = 0.002564847284773193
$ws.Range("M2").Value = 0.02094182787080249
$ws.Range("N2").Value = 0.03318585403457712
$ws.Range("O2").Value = 0.07551591960538651
$ws.Range("Q2").Value = 0.00232364718244743
$ws.Range("R2").Value = 0.00135184277276667
$ws.Range("S2").Value = 0.04781751904647121
$ws.Range("T2").Value = 0.009471119733151075
$ws.Range("U2").Value = 0.007896907875347242
$ws.Range("V2").Value = 0.01482658198043799
$ws.Range("W2").Value = 0.006229848235954951
$ws.Range("X2").Value = 0.03344313585766082
$ws.Range("Z2").Value = 0.03825674859410075
$ws.Range("AA3").Value = 0.01887965745992456
$ws.Range("AD3").Value = 0.02020081943296893
$ws.Range("AF3").Value = 0.02706307369575596
$ws.Range("AI3").Value = 0.0280320910704844
$ws.Range("D3").Value = 0.1924170420662693
$ws.Range("F3").Value = 0.4366347816105404
$ws.Range("G3").Value = 0.005627242830179106
$ws.Range("H3").Value = 0.02883431227065181
$ws.Range("J3").Value = 0.0003362608765465273
$ws.Range("M3").Value = 0.01877325983677277
$ws.Range("O3").Value = 0.1017730962087946
$ws.Range("R3").Value = 0.004187449504442063
$ws.Range("S3").Value = 0.01221814033626422
$ws.Range("T3").Value = 0.02667425976951174
$ws.Range("W3").Value = 0.01595822761459552
$ws.Range("X3").Value = 0.01536062667281516
$ws.Range("Y3").Value = 0.0003790855352055211
$ws.Range("Z3").Value = 0.04665057320827763
$ws.Range("AA4").Value = 0.01996676605840529
$ws.Range("AD4").Value = 0.02221759987258778
$ws.Range("AF4").Value = 0.01709817182483996
$ws.Range("AI4").Value = 0.01961942164720035
$ws.Range("D4").Value = 0.1059537519554112
$ws.Range("E4").Value = 0.04472007515110583
$ws.Range("F4").Value = 0.4048766031435542
$ws.Range("G4").Value = 0.05076488104615923
$ws.Range("H4").Value = 0.04050439371027997
$ws.Range("I4").Value = 0.01145499902702337
$ws.Range("L4").Value = 0.005944314844051236
$ws.Range("M4").Value = 0.00564112218012854
$ws.Range("O4").Value = 0.07811408909051817
$ws.Range("P4").Value = 0.004808440092791501
$ws.Range("R4").Value = 0.01055309595032341
$ws.Range("S4").Value = 0.007928092538535551
$ws.Range("T4").Value = 0.04328676465089256
$ws.Range("U4").Value = 0.001621783703160222
$ws.Range("W4").Value = 0.03674005073055029
$ws.Range("X4").Value = 0.01050535288958179
$ws.Range("Y4").Value = 0.009891634391994001
$ws.Range("Z4").Value = 0.04778859550090565
$ws.Range("AA5").Value = 0.009942631553762757
$ws.Range("AB5").Value = 0.01922845299703653
$ws.Range("AD5").Value = 0.01363600956642647
$ws.Range("AE5").Value = 0.009596973651371422
$ws.Range("AF5").Value = 0.02247865974913569
$ws.Range("AH5").Value = 0.0003795896105270282
$ws.Range("AI5").Value = 0.02608913051352927
$ws.Range("D5").Value = 0.2422211585652941
$ws.Range("E5").Value = 0.107832746582354
$ws.Range("F5").Value = 0.2249359574385186
$ws.Range("G5").Value = 0.006679026200989892
$ws.Range("H5").Value = 0.02519742967019477
$ws.Range("J5").Value = 0.01509187018188523
$ws.Range("K5").Value = 0.01700446836819745
$ws.Range("M5").Value = 0.02346579435141502
$ws.Range("N5").Value = 0.03037796710525244
$ws.Range("O5").Value = 0.07983715242072018
$ws.Range("S5").Value = 0.02056652975216587
$ws.Range("T5").Value = 0.007421861001979395
$ws.Range("V5").Value = 0.01269683711893177
$ws.Range("W5").Value = 0.003213977282477316
$ws.Range("X5").Value = 0.02838930981329043
$ws.Range("Z5").Value = 0.05371646650454448
$ws.Range("AA6").Value = 0.009103919018376637
$ws.Range("AB6").Value = 0.0135764556109032
$ws.Range("AD6").Value = 0.01675709219087264
$ws.Range("AE6").Value = 0.004184283519911789
$ws.Range("AF6").Value = 0.01775190882686376
$ws.Range("AI6").Value = 0.02467337056182232
$ws.Range("D6").Value = 0.2378787880602266
$ws.Range("E6").Value = 0.05340873513557569
$ws.Range("F6").Value = 0.3036246589961589
$ws.Range("G6").Value = 0.0008754550054492817
$ws.Range("H6").Value = 0.04185678710085436
$ws.Range("J6").Value = 0.0100937721874494
$ws.Range("K6").Value = 0.001019892967449335
$ws.Range("M6").Value = 0.02286840804307506
$ws.Range("N6").Value = 0.01562042295579223
$ws.Range("O6").Value = 0.09378277254539627
$ws.Range("S6").Value = 0.0144181984972618
$ws.Range("T6").Value = 0.00779999338368171
$ws.Range("V6").Value = 0.009674435227779193
$ws.Range("W6").Value = 0.005760935157820869
$ws.Range("X6").Value = 0.03178347320945011
$ws.Range("Z6").Value = 0.06348624179782869
$ws.Range("AA7").Value = 0.040010234451952
$ws.Range("AB7").Value = 0.008798933048156521
$ws.Range("AD7").Value = 0.02991001589362737
$ws.Range("AE7").Value = 0.005376797620421927
$ws.Range("AF7").Value = 0.07378309325176312
$ws.Range("AI7").Value = 0.01962630727994714
$ws.Range("D7").Value = 0.2716031974330277
$ws.Range("E7").Value = 0.08436418054987638
$ws.Range("F7").Value = 0.150693920716901
$ws.Range("I7").Value = 0.0002165397846571985
$ws.Range("J7").Value = 0.01801076621131503
$ws.Range("K7").Value = 0.0697893797109362
$ws.Range("M7").Value = 0.02598296444637585
$ws.Range("N7").Value = 0.03059550395212213
$ws.Range("O7").Value = 0.109793178662871
$ws.Range("Q7").Value = 0.0003371955493655521
$ws.Range("S7").Value = 0.01614524265429765
$ws.Range("V7").Value = 0.002669548798533916
$ws.Range("W7").Value = 0.0018395591131222
$ws.Range("X7").Value = 0.00106642905698535
$ws.Range("Y7").Value = 0.007645732192618827
$ws.Range("Z7").Value = 0.03174127962112586
$ws.Range("AA8").Value = 0.03866013691633959
$ws.Range("AB8").Value = 0.004990631215198611
$ws.Range("AC8").Value = 0.0005800333182918083
$ws.Range("AD8").Value = 0.03657459806350123
$ws.Range("AE8").Value = 0.003413904941063519
$ws.Range("AF8").Value = 0.06691695980968149
$ws.Range("AG8").Value = 0.0004451883569815222
$ws.Range("AI8").Value = 0.01936954081286194
$ws.Range("D8").Value = 0.2106756628772977
$ws.Range("E8").Value = 0.03469687696533146
$ws.Range("F8").Value = 0.2360681349489719
$ws.Range("G8").Value = 0.003800317309409391
$ws.Range("I8").Value = 0.003035920536644409
$ws.Range("J8").Value = 0.02390358164647598
$ws.Range("K8").Value = 0.06408383744562071
$ws.Range("L8").Value = 0.004391950257494481
$ws.Range("M8").Value = 0.05054425254308053
$ws.Range("N8").Value = 0.02027592505149497
$ws.Range("O8").Value = 0.1127124781978672
$ws.Range("Q8").Value = 0.00283306957817605
$ws.Range("S8").Value = 0.01719811823156463
$ws.Range("V8").Value = 0.004344003853789682
$ws.Range("W8").Value = 0.006076706673751485
$ws.Range("X8").Value = 0.003939752181172994
$ws.Range("Y8").Value = 0.005960213969700644
$ws.Range("Z8").Value = 0.02450820429823583
$ws.Range("AA9").Value = 0.0369018871166114
$ws.Range("AB9").Value = 0.001622885131564411
$ws.Range("AD9").Value = 0.04329167828616117
$ws.Range("AE9").Value = 0.001390655562317052
$ws.Range("AF9").Value = 0.06074500980720678
$ws.Range("AI9").Value = 0.02113366705817745
$ws.Range("D9").Value = 0.198631655127862
$ws.Range("E9").Value = 0.006833869992026506
$ws.Range("F9").Value = 0.2815999634778395
$ws.Range("J9").Value = 0.02331716145200837
$ws.Range("K9").Value = 0.05197631263843989
$ws.Range("L9").Value = 0.007968073198572604
$ws.Range("M9").Value = 0.04558654761938071
$ws.Range("N9").Value = 0.01967325024503735
$ws.Range("O9").Value = 0.1204899620933816
$ws.Range("Q9").Value = 0.005182404025865934
$ws.Range("S9").Value = 0.02261084220675576
$ws.Range("V9").Value = 0.005319254222642862
$ws.Range("W9").Value = 0.008383632020752322
$ws.Range("X9").Value = 0.001560146526170076
$ws.Range("Y9").Value = 0.004483142798021395
$ws.Range("Z9").Value = 0.03129799939320484
$ws.Range("AA10").Value = 0.04740785186025942
$ws.Range("AC10").Value = 0.004637889309083848
$ws.Range("AD10").Value = 0.04283231605633148
$ws.Range("AF10").Value = 0.06763139143192931
$ws.Range("AG10").Value = 0.003815424603664537
$ws.Range("AI10").Value = 0.01724405951909062
$ws.Range("D10").Value = 0.1930670815396631
$ws.Range("E10").Value = 0.001621225922547748
$ws.Range("F10").Value = 0.2994769303481208
$ws.Range("J10").Value = 0.02813758103578399
$ws.Range("K10").Value = 0.03548030358633512
$ws.Range("L10").Value = 0.02254991991892535
$ws.Range("M10").Value = 0.0468463126881376
$ws.Range("N10").Value = 0.006749988257526516
$ws.Range("O10").Value = 0.1232568435890316
$ws.Range("P10").Value = 0.00356720438500195
$ws.Range("S10").Value = 0.0110669825140757
$ws.Range("V10").Value = 0.002870534317371753
$ws.Range("W10").Value = 0.01074673216191237
$ws.Range("Y10").Value = 0.008280579994891233
$ws.Range("Z10").Value = 0.02271284696031603
$ws.Range("AA11").Value = 0.02453236838593181
$ws.Range("AB11").Value = 0.01311469163808167
$ws.Range("AD11").Value = 0.02660656882413126
$ws.Range("AE11").Value = 0.01229213785047731
$ws.Range("AF11").Value = 0.06126816266718416
$ws.Range("AI11").Value = 0.02708514738698229
$ws.Range("D11").Value = 0.2587224494303021
$ws.Range("E11").Value = 0.08799611988913193
$ws.Range("F11").Value = 0.1704956685356178
$ws.Range("G11").Value = 0.003990599138047928
$ws.Range("I11").Value = 0.005462754418352743
$ws.Range("J11").Value = 0.01626144534254061
$ws.Range("K11").Value = 0.06348837691004856
$ws.Range("M11").Value = 0.02621602469088981
$ws.Range("N11").Value = 0.02957056353685756
$ws.Range("O11").Value = 0.08882270371679911
$ws.Range("Q11").Value = 0.001977596387755333
$ws.Range("S11").Value = 0.03220929943531362
$ws.Range("V11").Value = 0.008977048175553386
$ws.Range("W11").Value = 0.002780588387495766
$ws.Range("X11").Value = 0.004701857982630716
$ws.Range("Y11").Value = 0.0003331439933172122
$ws.Range("Z11").Value = 0.0330946832765572

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AA2").Value = 0.9268307318862259
$ws.Range("AB2").Value = 0.9287782291983947
$ws.Range("AC2").Value = 0.9287782291983947
$ws.Range("AD2").Value = 0.9408928498922438
$ws.Range("AE2").Value = 0.9451271060516597
$ws.Range("AF2").Value = 0.9703178946548805
$ws.Range("AG2").Value = 0.9703178946548805
$ws.Range("AH2").Value = 0.9707318356762652
$ws.Range("D2").Value = 0.2022278673415956
$ws.Range("E2").Value = 0.2097153487035711
$ws.Range("F2").Value = 0.5980591220809204
$ws.Range("G2").Value = 0.5980591220809204
$ws.Range("H2").Value = 0.6280117097002617
$ws.Range("I2").Value = 0.6280117097002617
$ws.Range("J2").Value = 0.6305765569850349
$ws.Range("K2").Value = 0.6305765569850349
$ws.Range("L2").Value = 0.6305765569850349
$ws.Range("M2").Value = 0.6515183848558374
$ws.Range("N2").Value = 0.6847042388904144
$ws.Range("O2").Value = 0.7602201584958009
$ws.Range("P2").Value = 0.7602201584958009
$ws.Range("Q2").Value = 0.7625438056782483
$ws.Range("R2").Value = 0.763895648451015
$ws.Range("S2").Value = 0.8117131674974862
$ws.Range("T2").Value = 0.8211842872306373
$ws.Range("U2").Value = 0.8290811951059845
$ws.Range("V2").Value = 0.8439077770864225
$ws.Range("W2").Value = 0.8501376253223775
$ws.Range("X2").Value = 0.8835807611800383
$ws.Range("Y2").Value = 0.8835807611800383
$ws.Range("Z2").Value = 0.921837509774139
$ws.Range("AA3").Value = 0.9247040158007909
$ws.Range("AB3").Value = 0.9247040158007909
$ws.Range("AC3").Value = 0.9247040158007909
$ws.Range("AD3").Value = 0.9449048352337598
$ws.Range("AE3").Value = 0.9449048352337598
$ws.Range("AF3").Value = 0.9719679089295158
$ws.Range("AG3").Value = 0.9719679089295158
$ws.Range("AH3").Value = 0.9719679089295158
$ws.Range("D3").Value = 0.1924170420662693
$ws.Range("E3").Value = 0.1924170420662693
$ws.Range("F3").Value = 0.6290518236768097
$ws.Range("G3").Value = 0.6346790665069888
$ws.Range("H3").Value = 0.6635133787776406
$ws.Range("I3").Value = 0.6635133787776406
$ws.Range("J3").Value = 0.6638496396541871
$ws.Range("K3").Value = 0.6638496396541871
$ws.Range("L3").Value = 0.6638496396541871
$ws.Range("M3").Value = 0.6826228994909599
$ws.Range("N3").Value = 0.6826228994909599
$ws.Range("O3").Value = 0.7843959956997545
$ws.Range("P3").Value = 0.7843959956997545
$ws.Range("Q3").Value = 0.7843959956997545
$ws.Range("R3").Value = 0.7885834452041965
$ws.Range("S3").Value = 0.8008015855404608
$ws.Range("T3").Value = 0.8274758453099725
$ws.Range("U3").Value = 0.8274758453099725
$ws.Range("V3").Value = 0.8274758453099725
$ws.Range("W3").Value = 0.8434340729245681
$ws.Range("X3").Value = 0.8587946995973832
$ws.Range("Y3").Value = 0.8591737851325887
$ws.Range("Z3").Value = 0.9058243583408663
$ws.Range("AA4").Value = 0.9410648066553721
$ws.Range("AB4").Value = 0.9410648066553721
$ws.Range("AC4").Value = 0.9410648066553721
$ws.Range("AD4").Value = 0.9632824065279598
$ws.Range("AE4").Value = 0.9632824065279598
$ws.Range("AF4").Value = 0.9803805783527998
$ws.Range("AG4").Value = 0.9803805783527998
$ws.Range("AH4").Value = 0.9803805783527998
$ws.Range("D4").Value = 0.1059537519554112
$ws.Range("E4").Value = 0.150673827106517
$ws.Range("F4").Value = 0.5555504302500712
$ws.Range("G4").Value = 0.6063153112962305
$ws.Range("H4").Value = 0.6468197050065104
$ws.Range("I4").Value = 0.6582747040335338
$ws.Range("J4").Value = 0.6582747040335338
$ws.Range("K4").Value = 0.6582747040335338
$ws.Range("L4").Value = 0.6642190188775851
$ws.Range("M4").Value = 0.6698601410577136
$ws.Range("N4").Value = 0.6698601410577136
$ws.Range("O4").Value = 0.7479742301482318
$ws.Range("P4").Value = 0.7527826702410233
$ws.Range("Q4").Value = 0.7527826702410233
$ws.Range("R4").Value = 0.7633357661913467
$ws.Range("S4").Value = 0.7712638587298822
$ws.Range("T4").Value = 0.8145506233807748
$ws.Range("U4").Value = 0.816172407083935
$ws.Range("V4").Value = 0.816172407083935
$ws.Range("W4").Value = 0.8529124578144853
$ws.Range("X4").Value = 0.8634178107040671
$ws.Range("Y4").Value = 0.8733094450960611
$ws.Range("Z4").Value = 0.9210980405969668
$ws.Range("AA5").Value = 0.9085911839119738
$ws.Range("AB5").Value = 0.9278196369090104
$ws.Range("AC5").Value = 0.9278196369090104
$ws.Range("AD5").Value = 0.9414556464754369
$ws.Range("AE5").Value = 0.9510526201268084
$ws.Range("AF5").Value = 0.973531279875944
$ws.Range("AG5").Value = 0.973531279875944
$ws.Range("AH5").Value = 0.973910869486471
$ws.Range("D5").Value = 0.2422211585652941
$ws.Range("E5").Value = 0.3500539051476481
$ws.Range("F5").Value = 0.5749898625861667
$ws.Range("G5").Value = 0.5816688887871566
$ws.Range("H5").Value = 0.6068663184573514
$ws.Range("I5").Value = 0.6068663184573514
$ws.Range("J5").Value = 0.6219581886392366
$ws.Range("K5").Value = 0.638962657007434
$ws.Range("L5").Value = 0.638962657007434
$ws.Range("M5").Value = 0.6624284513588491
$ws.Range("N5").Value = 0.6928064184641015
$ws.Range("O5").Value = 0.7726435708848217
$ws.Range("P5").Value = 0.7726435708848217
$ws.Range("Q5").Value = 0.7726435708848217
$ws.Range("R5").Value = 0.7726435708848217
$ws.Range("S5").Value = 0.7932101006369876
$ws.Range("T5").Value = 0.800631961638967
$ws.Range("U5").Value = 0.800631961638967
$ws.Range("V5").Value = 0.8133287987578988
$ws.Range("W5").Value = 0.8165427760403761
$ws.Range("X5").Value = 0.8449320858536665
$ws.Range("Y5").Value = 0.8449320858536665
$ws.Range("Z5").Value = 0.8986485523582111
$ws.Range("AA6").Value = 0.9230568892896261
$ws.Range("AB6").Value = 0.9366333449005293
$ws.Range("AC6").Value = 0.9366333449005293
$ws.Range("AD6").Value = 0.9533904370914019
$ws.Range("AE6").Value = 0.9575747206113137
$ws.Range("AF6").Value = 0.9753266294381774
$ws.Range("AG6").Value = 0.9753266294381774
$ws.Range("AH6").Value = 0.9753266294381774
$ws.Range("AI6").Value = 0.9999999999999997
$ws.Range("D6").Value = 0.2378787880602266
$ws.Range("E6").Value = 0.2912875231958023
$ws.Range("F6").Value = 0.5949121821919612
$ws.Range("G6").Value = 0.5957876371974105
$ws.Range("H6").Value = 0.6376444242982648
$ws.Range("I6").Value = 0.6376444242982648
$ws.Range("J6").Value = 0.6477381964857143
$ws.Range("K6").Value = 0.6487580894531636
$ws.Range("L6").Value = 0.6487580894531636
$ws.Range("M6").Value = 0.6716264974962387
$ws.Range("N6").Value = 0.687246920452031
$ws.Range("O6").Value = 0.7810296929974272
$ws.Range("P6").Value = 0.7810296929974272
$ws.Range("Q6").Value = 0.7810296929974272
$ws.Range("R6").Value = 0.7810296929974272
$ws.Range("S6").Value = 0.795447891494689
$ws.Range("T6").Value = 0.8032478848783707
$ws.Range("U6").Value = 0.8032478848783707
$ws.Range("V6").Value = 0.81292232010615
$ws.Range("W6").Value = 0.8186832552639708
$ws.Range("X6").Value = 0.8504667284734209
$ws.Range("Y6").Value = 0.8504667284734209
$ws.Range("Z6").Value = 0.9139529702712496
$ws.Range("AA7").Value = 0.8625048529060839
$ws.Range("AB7").Value = 0.8713037859542404
$ws.Range("AC7").Value = 0.8713037859542404
$ws.Range("AD7").Value = 0.9012138018478678
$ws.Range("AE7").Value = 0.9065905994682898
$ws.Range("AF7").Value = 0.9803736927200529
$ws.Range("AG7").Value = 0.9803736927200529
$ws.Range("AH7").Value = 0.9803736927200529
$ws.Range("D7").Value = 0.2716031974330277
$ws.Range("E7").Value = 0.3559673779829041
$ws.Range("F7").Value = 0.5066612986998051
$ws.Range("G7").Value = 0.5066612986998051
$ws.Range("H7").Value = 0.5066612986998051
$ws.Range("I7").Value = 0.5068778384844623
$ws.Range("J7").Value = 0.5248886046957774
$ws.Range("K7").Value = 0.5946779844067136
$ws.Range("L7").Value = 0.5946779844067136
$ws.Range("M7").Value = 0.6206609488530894
$ws.Range("N7").Value = 0.6512564528052116
$ws.Range("O7").Value = 0.7610496314680827
$ws.Range("P7").Value = 0.7610496314680827
$ws.Range("Q7").Value = 0.7613868270174482
$ws.Range("R7").Value = 0.7613868270174482
$ws.Range("S7").Value = 0.7775320696717458
$ws.Range("T7").Value = 0.7775320696717458
$ws.Range("U7").Value = 0.7775320696717458
$ws.Range("V7").Value = 0.7802016184702797
$ws.Range("W7").Value = 0.7820411775834019
$ws.Range("X7").Value = 0.7831076066403873
$ws.Range("Y7").Value = 0.7907533388330061
$ws.Range("Z7").Value = 0.822494618454132
$ws.Range("AA8").Value = 0.8677091434824197
$ws.Range("AB8").Value = 0.8726997746976183
$ws.Range("AC8").Value = 0.8732798080159101
$ws.Range("AD8").Value = 0.9098544060794114
$ws.Range("AE8").Value = 0.9132683110204749
$ws.Range("AF8").Value = 0.9801852708301564
$ws.Range("AG8").Value = 0.9806304591871379
$ws.Range("AH8").Value = 0.9806304591871379
$ws.Range("AI8").Value = 0.9999999999999998
$ws.Range("D8").Value = 0.2106756628772977
$ws.Range("E8").Value = 0.2453725398426291
$ws.Range("F8").Value = 0.481440674791601
$ws.Range("G8").Value = 0.4852409921010104
$ws.Range("H8").Value = 0.4852409921010104
$ws.Range("I8").Value = 0.4882769126376548
$ws.Range("J8").Value = 0.5121804942841308
$ws.Range("K8").Value = 0.5762643317297516
$ws.Range("L8").Value = 0.5806562819872461
$ws.Range("M8").Value = 0.6312005345303265
$ws.Range("N8").Value = 0.6514764595818215
$ws.Range("O8").Value = 0.7641889377796888
$ws.Range("P8").Value = 0.7641889377796888
$ws.Range("Q8").Value = 0.7670220073578649
$ws.Range("R8").Value = 0.7670220073578649
$ws.Range("S8").Value = 0.7842201255894294
$ws.Range("T8").Value = 0.7842201255894294
$ws.Range("U8").Value = 0.7842201255894294
$ws.Range("V8").Value = 0.7885641294432191
$ws.Range("W8").Value = 0.7946408361169706
$ws.Range("X8").Value = 0.7985805882981436
$ws.Range("Y8").Value = 0.8045408022678442
$ws.Range("Z8").Value = 0.82904900656608
$ws.Range("AA9").Value = 0.8718161041545732
$ws.Range("AB9").Value = 0.8734389892861376
$ws.Range("AC9").Value = 0.8734389892861376
$ws.Range("AD9").Value = 0.9167306675722988
$ws.Range("AE9").Value = 0.9181213231346159
$ws.Range("AF9").Value = 0.9788663329418227
$ws.Range("AG9").Value = 0.9788663329418227
$ws.Range("AH9").Value = 0.9788663329418227
$ws.Range("AI9").Value = 1
$ws.Range("D9").Value = 0.198631655127862
$ws.Range("E9").Value = 0.2054655251198885
$ws.Range("F9").Value = 0.487065488597728
$ws.Range("G9").Value = 0.487065488597728
$ws.Range("H9").Value = 0.487065488597728
$ws.Range("I9").Value = 0.487065488597728
$ws.Range("J9").Value = 0.5103826500497364
$ws.Range("K9").Value = 0.5623589626881763
$ws.Range("L9").Value = 0.5703270358867489
$ws.Range("M9").Value = 0.6159135835061297
$ws.Range("N9").Value = 0.635586833751167
$ws.Range("O9").Value = 0.7560767958445486
$ws.Range("P9").Value = 0.7560767958445486
$ws.Range("Q9").Value = 0.7612591998704145
$ws.Range("R9").Value = 0.7612591998704145
$ws.Range("S9").Value = 0.7838700420771703
$ws.Range("T9").Value = 0.7838700420771703
$ws.Range("U9").Value = 0.7838700420771703
$ws.Range("V9").Value = 0.7891892962998132
$ws.Range("W9").Value = 0.7975729283205655
$ws.Range("X9").Value = 0.7991330748467356
$ws.Range("Y9").Value = 0.803616217644757
$ws.Range("Z9").Value = 0.8349142170379619
$ws.Range("AA10").Value = 0.8638389190799001
$ws.Range("AB10").Value = 0.8638389190799001
$ws.Range("AC10").Value = 0.8684768083889839
$ws.Range("AD10").Value = 0.9113091244453154
$ws.Range("AE10").Value = 0.9113091244453154
$ws.Range("AF10").Value = 0.9789405158772447
$ws.Range("AG10").Value = 0.9827559404809093
$ws.Range("AH10").Value = 0.9827559404809093
$ws.Range("D10").Value = 0.1930670815396631
$ws.Range("E10").Value = 0.1946883074622109
$ws.Range("F10").Value = 0.4941652378103317
$ws.Range("G10").Value = 0.4941652378103317
$ws.Range("H10").Value = 0.4941652378103317
$ws.Range("I10").Value = 0.4941652378103317
$ws.Range("J10").Value = 0.5223028188461156
$ws.Range("K10").Value = 0.5577831224324508
$ws.Range("L10").Value = 0.5803330423513761
$ws.Range("M10").Value = 0.6271793550395137
$ws.Range("N10").Value = 0.6339293432970402
$ws.Range("O10").Value = 0.7571861868860718
$ws.Range("P10").Value = 0.7607533912710737
$ws.Range("Q10").Value = 0.7607533912710737
$ws.Range("R10").Value = 0.7607533912710737
$ws.Range("S10").Value = 0.7718203737851493
$ws.Range("T10").Value = 0.7718203737851493
$ws.Range("U10").Value = 0.7718203737851493
$ws.Range("V10").Value = 0.7746909081025211
$ws.Range("W10").Value = 0.7854376402644334
$ws.Range("X10").Value = 0.7854376402644334
$ws.Range("Y10").Value = 0.7937182202593246
$ws.Range("Z10").Value = 0.8164310672196406
$ws.Range("AA11").Value = 0.8596332916331428
$ws.Range("AB11").Value = 0.8727479832712245
$ws.Range("AC11").Value = 0.8727479832712245
$ws.Range("AD11").Value = 0.8993545520953558
$ws.Range("AE11").Value = 0.9116466899458331
$ws.Range("AF11").Value = 0.9729148526130172
$ws.Range("AG11").Value = 0.9729148526130172
$ws.Range("AH11").Value = 0.9729148526130172
$ws.Range("D11").Value = 0.2587224494303021
$ws.Range("E11").Value = 0.346718569319434
$ws.Range("F11").Value = 0.5172142378550517
$ws.Range("G11").Value = 0.5212048369930996
$ws.Range("H11").Value = 0.5212048369930996
$ws.Range("I11").Value = 0.5266675914114524
$ws.Range("J11").Value = 0.542929036753993
$ws.Range("K11").Value = 0.6064174136640416
$ws.Range("L11").Value = 0.6064174136640416
$ws.Range("M11").Value = 0.6326334383549314
$ws.Range("N11").Value = 0.662204001891789
$ws.Range("O11").Value = 0.7510267056085881
$ws.Range("P11").Value = 0.7510267056085881
$ws.Range("Q11").Value = 0.7530043019963434
$ws.Range("R11").Value = 0.7530043019963434
$ws.Range("S11").Value = 0.785213601431657
$ws.Range("T11").Value = 0.785213601431657
$ws.Range("U11").Value = 0.785213601431657
$ws.Range("V11").Value = 0.7941906496072103
$ws.Range("W11").Value = 0.796971237994706
$ws.Range("X11").Value = 0.8016730959773367
$ws.Range("Y11").Value = 0.8020062399706539
$ws.Range("Z11").Value = 0.8351009232472111

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.5980591220809204
$ws.Range("F3").Value = 0.6290518236768097
$ws.Range("F4").Value = 0.5555504302500712
$ws.Range("F5").Value = 0.5749898625861667
$ws.Range("F6").Value = 0.5949121821919612
$ws.Range("F7").Value = 0.5066612986998051
$ws.Range("D8").Value = 9
$ws.Range("F8").Value = 0.5121804942841308
$ws.Range("G8").Value = 8
$ws.Range("D9").Value = 9
$ws.Range("F9").Value = 0.5103826500497364
$ws.Range("G9").Value = 8
$ws.Range("D10").Value = 9
$ws.Range("F10").Value = 0.5223028188461156
$ws.Range("G10").Value = 8
$ws.Range("F11").Value = 0.5172142378550517

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 14
$ws.Range("F2").Value = 0.7602201584958009
$ws.Range("G2").Value = 13
$ws.Range("D3").Value = 14
$ws.Range("F3").Value = 0.7843959956997545
$ws.Range("G3").Value = 13
$ws.Range("D4").Value = 14
$ws.Range("F4").Value = 0.7479742301482318
$ws.Range("G4").Value = 13
$ws.Range("D5").Value = 14
$ws.Range("F5").Value = 0.7726435708848217
$ws.Range("G5").Value = 13
$ws.Range("D6").Value = 14
$ws.Range("F6").Value = 0.7810296929974272
$ws.Range("G6").Value = 13
$ws.Range("F7").Value = 0.7610496314680827
$ws.Range("F8").Value = 0.7641889377796888
$ws.Range("F9").Value = 0.7560767958445486
$ws.Range("F10").Value = 0.7571861868860718
$ws.Range("D11").Value = 14
$ws.Range("F11").Value = 0.7510267056085881
$ws.Range("G11").Value = 13

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 18
$ws.Range("F2").Value = 0.8117131674974862
$ws.Range("G2").Value = 17
$ws.Range("D3").Value = 18
$ws.Range("F3").Value = 0.8008015855404608
$ws.Range("G3").Value = 17
$ws.Range("D4").Value = 19
$ws.Range("F4").Value = 0.8145506233807748
$ws.Range("G4").Value = 18
$ws.Range("D5").Value = 19
$ws.Range("F5").Value = 0.800631961638967
$ws.Range("G5").Value = 18
$ws.Range("D6").Value = 19
$ws.Range("F6").Value = 0.8032478848783707
$ws.Range("G6").Value = 18
$ws.Range("D7").Value = 25
$ws.Range("F7").Value = 0.822494618454132
$ws.Range("G7").Value = 24
$ws.Range("D8").Value = 24
$ws.Range("F8").Value = 0.8045408022678442
$ws.Range("G8").Value = 23
$ws.Range("D9").Value = 24
$ws.Range("F9").Value = 0.803616217644757
$ws.Range("G9").Value = 23
$ws.Range("D10").Value = 25
$ws.Range("F10").Value = 0.8164310672196406
$ws.Range("G10").Value = 24
$ws.Range("D11").Value = 23
$ws.Range("F11").Value = 0.8016730959773367
$ws.Range("G11").Value = 22

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 25
$ws.Range("F2").Value = 0.921837509774139
$ws.Range("G2").Value = 24
$ws.Range("F3").Value = 0.9058243583408663
$ws.Range("F4").Value = 0.9210980405969668
$ws.Range("D5").Value = 26
$ws.Range("F5").Value = 0.9085911839119738
$ws.Range("G5").Value = 25
$ws.Range("F6").Value = 0.9139529702712496
$ws.Range("F7").Value = 0.9012138018478678
$ws.Range("F8").Value = 0.9098544060794114
$ws.Range("F9").Value = 0.9167306675722988
$ws.Range("F10").Value = 0.9113091244453154
$ws.Range("D11").Value = 30
$ws.Range("F11").Value = 0.9116466899458331
$ws.Range("G11").Value = 29
